# Apply updated cryptocurrency price/volume data to Sheet1
# (values that look like plain decimal numbers are forced to
#  Text format first so Excel does not auto-convert/round them)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '68.179.21'
$ws.Range('E2').Value = '  -3.61%  '
# Row 3
$ws.Range('D3').Value = '3.698.23'
$ws.Range('E3').Value = '  -4.31%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  +0.04%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.05'
$ws.Range('E5').Value = '  +0.93%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.29'
$ws.Range('E6').Value = '  +8.33%  '
# Row 7
$ws.Range('D7').Value = '3.692.21'
$ws.Range('E7').Value = '  -4.22%  '
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.626'
$ws.Range('E8').Value = '  -6.47%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.997'
$ws.Range('E9').Value = '  -0.38%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.717'
$ws.Range('E10').Value = '  -4.38%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.163'
$ws.Range('E11').Value = '  -6.80%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '56.08'
$ws.Range('E12').Value = '  +4.77%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000291'
$ws.Range('E13').Value = '  -9.53%  '
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.40'
$ws.Range('E14').Value = '  -7.97%  '
# Row 15
$ws.Range('D15').Value = '4.266.57'
$ws.Range('E15').Value = '  -4.92%  '
# Row 16
$ws.Range('D16').Value = '3.685.98'
$ws.Range('E16').Value = '  -4.31%  '
# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.42'
$ws.Range('E17').Value = '  -6.14%  '
# Row 18
$ws.Range('E18').Value = '  -2.21%  '
# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.82'
$ws.Range('E19').Value = '  -7.23%  '
# Row 20
$ws.Range('E20').Value = '  -6.68%  '
# Row 21
$ws.Range('D21').Value = '67.867.51'
$ws.Range('E21').Value = '  -3.90%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '410.26'
$ws.Range('E22').Value = '  -5.86%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.55'
$ws.Range('E23').Value = '  -3.77%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.60'
$ws.Range('E24').Value = '  -5.96%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.01'
$ws.Range('E25').Value = '  -7.99%  '
# Row 26
$ws.Range('B26').Value = 'RenderToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.06'
$ws.Range('E26').Value = '  -0.31%  '
# Row 27
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.82'
$ws.Range('E27').Value = '  -7.10%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.87'
$ws.Range('E28').Value = '  -4.34%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.06'
$ws.Range('E29').Value = '  +2.25%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.48'
$ws.Range('E30').Value = '  -7.87%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.74'
$ws.Range('E31').Value = '  -6.55%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.30'
$ws.Range('E32').Value = '  -8.87%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '12.51'
$ws.Range('E33').Value = '  -7.39%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.117'
$ws.Range('E34').Value = '  -6.28%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '43.55'
$ws.Range('E35').Value = '  -10.97%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '64.22'
$ws.Range('E36').Value = '  -8.31%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '590.81'
$ws.Range('E37').Value = '  -5.68%  '
# Row 38
$ws.Range('D38').Value = '0.0₃0880'
$ws.Range('E38').Value = '  -9.78%  '
# Row 39
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.401'
$ws.Range('E39').Value = '  -4.63%  '
# Row 40
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.10%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('E41').Value = '  -0.02%  '
# Row 42
$ws.Range('E42').Value = '  -5.02%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.78'
$ws.Range('E43').Value = '  +2.28%  '
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.00'
$ws.Range('E44').Value = '  -7.96%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0434'
$ws.Range('E45').Value = '  -7.43%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.87'
$ws.Range('E46').Value = '  -13.65%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.21'
$ws.Range('E47').Value = '  -8.70%  '
# Row 48
$ws.Range('E48').Value = '  -3.60%  '
# Row 49
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.134'
$ws.Range('E49').Value = '  -6.61%  '
# Row 50
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.766.11'
$ws.Range('E50').Value = '  -2.48%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.14'
$ws.Range('E51').Value = '  -4.39%  '
